$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.870.55"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.741.67"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.18"
$ws.Range("E5").Value = "  +3.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5152"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2742"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.02"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06128"
$ws.Range("D11").Value = "1.740.36"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07180"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.01"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6446"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.599"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.41"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9989"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9985"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "25.886.44"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.77"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006781"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "1.962.56"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.274"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.676"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.251"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.81"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.508"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.25"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.771"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "106.30"
$ws.Range("E30").Value = "  +4.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.027"
$ws.Range("E31").Value = "  +9.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08324"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.651"
$ws.Range("E33").Value = "  +4.34%  "
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.661"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9919"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6210"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.687"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01619"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.936"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9987"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.15"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3858"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7356"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.956"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1127"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.199"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05256"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.08"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.597"
$ws.Range("E51").Value = "  -0.06%  "

Write-Host "Applied cryptos update"
